# Rewrites "produit Expat" (drops the numeric index column, keeps Nom/Prix
# Produit/Prix Promo) and adds a brand-new "Produit Soumari" sheet scraped
# from a second e-commerce site (networking gear).
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "produit Expat": clear + rebuild with columns A=Nom, B=Prix Produit, C=Prix Promo
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("produit Expat")
$ws1.Cells.Clear()

$header1 = $ws1.Range("A1:C1")
$header1.Font.Bold = $true
$header1.Borders.LineStyle = 1
$header1.HorizontalAlignment = -4108   # xlCenter
$header1.VerticalAlignment = -4160     # xlTop
$ws1.Range("A1").Value = "Nom"
$ws1.Range("B1").Value = "Prix Produit"
$ws1.Range("C1").Value = "Prix Promo"

# Prices are stored as text ("0.00", "25000.00", ...) in the source data, so
# force column B to Text format before writing, otherwise Excel would coerce
# these numeric-looking strings into actual numbers.
$ws1.Range("B2:B12").NumberFormat = "@"

$ws1.Range("A2").Value = "Tenue de cérémonie à vider"
$ws1.Range("B2").Value = "0.00"
$ws1.Range("A3").Value = "Tenue de cérémonie à vider"
$ws1.Range("B3").Value = "0.00"
$ws1.Range("A4").Value = "ABAYAS"
$ws1.Range("B4").Value = "25000.00"
$ws1.Range("A5").Value = "Montres"
$ws1.Range("B5").Value = "27000.00"
$ws1.Range("A6").Value = "Lacoste pour homme"
$ws1.Range("B6").Value = "10000.00"
$ws1.Range("A7").Value = "Chaussure enfant homme et femme"
$ws1.Range("B7").Value = "0.00"
$ws1.Range("A8").Value = "Sac à main en wax"
$ws1.Range("B8").Value = "12000.00"
$ws1.Range("A9").Value = "Des thioups légers venant de la mauritanie"
$ws1.Range("B9").Value = "15000.00"
$ws1.Range("A10").Value = "Vêtements homme"
$ws1.Range("B10").Value = "4500.00"
$ws1.Range("A11").Value = "Sacoches disponibles neuves"
$ws1.Range("B11").Value = "3000.00"
$ws1.Range("A12").Value = "Sacoches"
$ws1.Range("B12").Value = "7500.00"

# ---------------------------------------------------------------------------
# 2. New sheet "Produit Soumari" (after ProduitJumia, i.e. last tab)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Produit Soumari"

$header3 = $ws3.Range("A1:B1")
$header3.Font.Bold = $true
$header3.Borders.LineStyle = 1
$header3.HorizontalAlignment = -4108   # xlCenter
$header3.VerticalAlignment = -4160     # xlTop
$ws3.Range("A1").Value = "Nom"
$ws3.Range("B1").Value = "Prix Produit"

$ws3.Range("B2:B49").NumberFormat = "@"

$ws3.Range("A2").Value = "Imprimante HP OFFICEJET 6950"
$ws3.Range("B2").Value = "72.000"
$ws3.Range("A3").Value = "Imprimante HP Laserjet PRO MFP M127FW"
$ws3.Range("B3").Value = "165.000"
$ws3.Range("A4").Value = "Imprimante HP DESKJET 3636"
$ws3.Range("B4").Value = "37.200"
$ws3.Range("A5").Value = "Imprimante HP DESKJET 2131"
$ws3.Range("B5").Value = "27.500"
$ws3.Range("A6").Value = "Imprimante HP MFP M479FDW"
$ws3.Range("B6").Value = "699.600"
$ws3.Range("A7").Value = "Imprimante HP Color MFP M176N"
$ws3.Range("B7").Value = "187.000"
$ws3.Range("A8").Value = "Imprimante DYMO Label Writer 450"
$ws3.Range("B8").Value = "60.000"
$ws3.Range("A9").Value = "Switch D LINK DGS-1024D"
$ws3.Range("B9").Value = "84.000"
$ws3.Range("A10").Value = "Switch D LINK DES-1024D"
$ws3.Range("B10").Value = "40.250"
$ws3.Range("A11").Value = "Switch D LINK DES-1005A/1005C"
$ws3.Range("B11").Value = "7.500"
$ws3.Range("A12").Value = "Switch UBIQUITI Lite USW-LITE-16-POE"
$ws3.Range("B12").Value = "243.000"
$ws3.Range("A13").Value = "Switch UBIQUITI UNIFI USW-PRO-48-POE"
$ws3.Range("B13").Value = "848.000"
$ws3.Range("A14").Value = "Switch UBIQUITI UNIFI USW-16-POE"
$ws3.Range("B14").Value = "342.400"
$ws3.Range("A15").Value = "Switch UBIQUITI UNIFI US-48-500W"
$ws3.Range("B15").Value = "636.000"
$ws3.Range("A16").Value = "Switch TP LINK TL-SG3452XP"
$ws3.Range("B16").Value = "529.650"
$ws3.Range("A17").Value = "Switch TP LINK TL-SG3452P"
$ws3.Range("B17").Value = "371.000"
$ws3.Range("A18").Value = "Switch TP LINK TL-SF1048 48 Ports"
$ws3.Range("B18").Value = "88.000"
$ws3.Range("A19").Value = "Switch TP LINK TL-SG1048 48 Ports"
$ws3.Range("B19").Value = "187.000"
$ws3.Range("A20").Value = "Switch TP LINK TL-SG1024D 24 Ports"
$ws3.Range("B20").Value = "78.000"
$ws3.Range("A21").Value = "Switch TP LINK TL-SG1005D 5 Ports"
$ws3.Range("B21").Value = "15.000"
$ws3.Range("A22").Value = "Switch TENDA TEG1016D"
$ws3.Range("B22").Value = "43.350"
$ws3.Range("A23").Value = "Switch TENDA SG108 8 Ports"
$ws3.Range("B23").Value = "25.000"
$ws3.Range("A24").Value = "Switch TENDA SG105 5 Ports"
$ws3.Range("B24").Value = "18.750"
$ws3.Range("A25").Value = "Switch TENDA TEG1105P-4-63W 4 Ports"
$ws3.Range("B25").Value = "43.350"
$ws3.Range("A26").Value = "Switch TENDA TEF1110P-8-102W 8 Ports"
$ws3.Range("B26").Value = "37.150"
$ws3.Range("A27").Value = "Switch D LINK DGS-1016D 16 Ports"
$ws3.Range("B27").Value = "49.500"
$ws3.Range("A28").Value = "Imprimante CANON MF752CDW"
$ws3.Range("B28").Value = "508.250"
$ws3.Range("A29").Value = "Imprimante HP Multifonction 137FNW"
$ws3.Range("B29").Value = "226.800"
$ws3.Range("A30").Value = "Imprimante EPSON ECOTANK L3251"
$ws3.Range("B30").Value = "176.000"
$ws3.Range("A31").Value = "Ensemble Clavier et Souris HP CS700"
$ws3.Range("B31").Value = "9.000"
$ws3.Range("A32").Value = "ROUTEUR D-LINK DIR-514"
$ws3.Range("B32").Value = "19.500"
$ws3.Range("A33").Value = "Modem Routeur TP LINK ARCHER D5"
$ws3.Range("B33").Value = "66.000"
$ws3.Range("A34").Value = "Modem Routeur TP LINK ARCHER MR200"
$ws3.Range("B34").Value = "88.000"
$ws3.Range("A35").Value = "Routeur Sans Fil TP LINK TL-MR3040"
$ws3.Range("B35").Value = "25.000"
$ws3.Range("A36").Value = "Boitier WiFi TP LINK DECO X60 (Pack de 3 Routeurs) sans fil"
$ws3.Range("B36").Value = "324.000"
$ws3.Range("A37").Value = "Routeur TP LINK ARCHER AX6000"
$ws3.Range("B37").Value = "259.200"
$ws3.Range("A38").Value = "Routeur WiFi TP LINK ARCHER C6"
$ws3.Range("B38").Value = "49.600"
$ws3.Range("A39").Value = "Routeur Bi-bande sans fil TP LINK ARCHER C58HP"
$ws3.Range("B39").Value = "66.000"
$ws3.Range("A40").Value = "Routeur WiFi TP LINK TL-WR840N"
$ws3.Range("B40").Value = "18.000"
$ws3.Range("A41").Value = "Routeur D LINK DSL-124"
$ws3.Range("B41").Value = "19.500"
$ws3.Range("A42").Value = "Modem Routeur TP LINK ARCHER VR400"
$ws3.Range("B42").Value = "52.650"
$ws3.Range("A43").Value = "Routeur Double Bande Sans Fil TP LINK MR3620"
$ws3.Range("B43").Value = "40.300"
$ws3.Range("A44").Value = "Routeur WiFi TP LINK ARCHER AX55/AX3000"
$ws3.Range("B44").Value = "88.000"
$ws3.Range("A45").Value = "Modem Routeur TENDA TDE-AC6"
$ws3.Range("B45").Value = "31.250"
$ws3.Range("A46").Value = "Modem Routeur TENDA TDE-AC10"
$ws3.Range("B46").Value = "37.200"
$ws3.Range("A47").Value = "Répéteur WiFi TP LINK RE300"
$ws3.Range("B47").Value = "43.400"
$ws3.Range("A48").Value = "Modem Routeur TENDA D301"
$ws3.Range("B48").Value = "19.500"
$ws3.Range("A49").Value = "Modem Routeur TENDA D305"
$ws3.Range("B49").Value = "25.000"

[void]$ws1.Range("A1").Select()
